# Apply the "application update with redux engine" edit to CheckList.xlsx
# - Remove the now-unused "SOSO" conditional-formatting rule / value
# - Flip most status cells in column C to NOK (keep the first few rows OK)
# - Append two new checklist rows: "language switch" and "theme switch"
# - Update the sheet selection to match the end of the new data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the standalone "SOSO" conditional formatting rule that applied to C2:C136.
#    (Range("C2").FormatConditions enumerates the 3 rules that apply to C2: NOK, OK, SOSO.)
$soSoRules = $ws.Range("C2").FormatConditions
for ($i = $soSoRules.Count; $i -ge 1; $i--) {
    $rule = $soSoRules.Item($i)
    if ($rule.Formula1 -like "*SOSO*") {
        $rule.Delete()
    }
}

# 2) Update column C values.
#    Rows 2-6 stay "OK" (already the case); row 2 is explicitly set for clarity.
$ws.Range("C2").Value = "OK"

#    Rows 7-55 (including the old "SOSO" row 35) become "NOK".
$ws.Range("C7:C55").Value = "NOK"

# 3) Append the two new checklist rows.
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = "language switch"
$ws.Range("C56").Value = "NOK"

$ws.Range("A57").Value = 56
$ws.Range("B57").Value = "theme switch"
$ws.Range("C57").Value = "NOK"

# 4) Update the active selection to highlight the newly added rows.
$ws.Range("A54:A57").Select()
